$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant for Paste Special "Formats only"
$xlPasteFormats = -4122

# --- 1. New TestCase4 row label (creates shared string "TestCase4") ---
$ws.Range("A5").Value = "TestCase4"

# --- 2. New header cells for the address/profile fields (I1:M1) ---
$ws.Range("I1").Value = "fname"
$ws.Range("J1").Value = "flocality"
$ws.Range("K1").Value = "faddress"
$ws.Range("L1").Value = "fcity"
$ws.Range("M1").Value = "fstate"

# --- 3. Fill in row 4 (D4:H4) with the same product/pincode data as row 3 ---
$ws.Range("D4").Value = "Bingo Mad Angles Achaari Masti Chips Chips"
$ws.Range("E4").Value = "Parrys White Label Sugar"
$ws.Range("F4").Value = "PARLE Wafers Cream and Onion"
$ws.Range("G4").Value = "'641402"
$ws.Range("H4").Value = "Parry's White Label Sugar"

# Copy the formatting (quote-prefixed text style) from row 3 onto the new row 4 cells
$ws.Range("D3:H3").Copy()
$ws.Range("D4:H4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# --- 4. Fill the rest of row 5 (B5:H5) to mirror rows 2-4 ---
# Add the hyperlink first (Excel forces the cell text to the display argument),
# then overwrite the cell value back to the properly-capitalised text.
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:vino@123", $null, $null, "vino@123")

$ws.Range("B5").Value = "'9952622682"
$ws.Range("C5").Value = "Vino@123"
$ws.Range("D5").Value = "Bingo Mad Angles Achaari Masti Chips Chips"
$ws.Range("E5").Value = "Parrys White Label Sugar"
$ws.Range("F5").Value = "PARLE Wafers Cream and Onion"
$ws.Range("G5").Value = "'641402"
$ws.Range("H5").Value = "Parry's White Label Sugar"

$ws.Range("B4:H4").Copy()
$ws.Range("B5:H5").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# --- 5. New profile/address columns for row 5 (I5:M5) ---
$ws.Range("I5").Value = "Vino"
$ws.Range("J5").Value = "coimbatore"
$ws.Range("M5").Value = "Tamil Nadu"
$ws.Range("K5").Value = "airforce station,sulur"
$ws.Range("L5").Value = "coimbatore"

$ws.Range("D4:H4").Copy()
$ws.Range("I5:M5").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# --- 6. Update the active selection like the saved workbook ---
$ws.Range("G20").Select()
